# Add finance LOC member
# Adds a new committee member row (Chantal Meré, BMS) with a finance (I)
# and loc_extended (F) flag, and marks F49 (the previous last row) as
# loc_extended too, matching the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting down from existing rows so the new cells pick up the
# --- same cell styles (s="7" for A:C, s="6" for F, s="9" for I) ---

# A50:C50 should look like A49:C49 (style index 7)
$ws.Range("A49:C49").Copy() | Out-Null
$ws.Range("A50:C50").PasteSpecial(-4122) | Out-Null

# F49 and F50 should look like F46/F47 (style index 6, centered "1" flag)
$ws.Range("F46").Copy() | Out-Null
$ws.Range("F49").PasteSpecial(-4122) | Out-Null
$ws.Range("F50").PasteSpecial(-4122) | Out-Null

# I50 should look like I49 (style index 9, centered "1" flag)
$ws.Range("I49").Copy() | Out-Null
$ws.Range("I50").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = 0

# --- Populate the actual values ---

# New committee member: Chantal Meré, institution BMS
$ws.Range("A50").Value = "Chantal"
$ws.Range("B50").Value = "Meré"
$ws.Range("C50").Value = "BMS"

# loc_extended flag for the newly added row and for the previous last row
$ws.Range("F49").Value = 1
$ws.Range("F50").Value = 1

# finance flag for the newly added row (existing I49 flag stays as-is)
$ws.Range("I49").Value = 1
$ws.Range("I50").Value = 1

# --- Update the sheet view so it scrolls back to the top and selects the
# --- new last cell, as in the authoritative workbook ---
$ws.Activate() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("J50").Select() | Out-Null
